# mother_baby_stock.xlsx — "final decided algo" update
#
# The one-year-return calculation (column M) was re-run with the final
# algorithm, which changed several return percentages. The colour-coding
# thresholds applied to the return (M) and live-price (D) columns were
# re-evaluated against the new numbers:
#   blue  (RGB 0,0,255)   -> mother_live_price / "neutral" one-year return
#   green (RGB 0,128,0)   -> positive / strong one-year return & live price
#   red   (RGB 255,0,0)   -> negative one-year return
#
# This script re-applies the recomputed values and the resulting colours
# via the Excel object model (Range.Value2 / Range.Font.Color), exactly as
# Excel itself would when a user (or a macro) edits the cells by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# OLE/VBA colour values (BGR-packed), equivalent to RGB(r,g,b) = r + g*256 + b*65536
$Blue  = 16711680   # RGB(0,0,255)
$Green = 32768       # RGB(0,128,0)
$Red   = 255         # RGB(255,0,0)

# --- Row 2 : Central Bank of India -----------------------------------
$ws.Range("M2").Value2 = 17.56
$ws.Range("M2").Font.Color = $Blue

# --- Row 3 : Equitas Small Finance Bank Ltd. --------------------------
$ws.Range("M3").Value2 = -33.83
$ws.Range("M3").Font.Color = $Red

# --- Row 4 : Jyoti CNC Automation Ltd. --------------------------------
$ws.Range("D4").Font.Color = $Green
$ws.Range("M4").Font.Color = $Green
# M4 value (156.56) is unchanged

# --- Row 5 : AIA Engineering Ltd. --------------------------------------
$ws.Range("D5").Font.Color = $Green
$ws.Range("M5").Value2 = 13.42
$ws.Range("M5").Font.Color = $Blue

# --- Row 6 : FSN E-Commerce Ventures Ltd. ------------------------------
$ws.Range("D6").Font.Color = $Green
$ws.Range("M6").Value2 = 7.13
$ws.Range("M6").Font.Color = $Blue

# --- Row 7 : KPIT Technologies Ltd. ------------------------------------
$ws.Range("D7").Font.Color = $Green
$ws.Range("M7").Value2 = 20.03
$ws.Range("M7").Font.Color = $Green

# --- Row 8 : MphasiS Ltd. -----------------------------------------------
$ws.Range("D8").Font.Color = $Green
$ws.Range("M8").Value2 = 15.56
$ws.Range("M8").Font.Color = $Blue
